$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "PDB molecule" column (column D) entirely; this shifts
# the "PDB filename" and "Is model" columns left by one.
$ws.Columns("D").Delete()

# After deleting the column, Excel leaves the selection on the column
# that now occupies D (previously E).
$ws.Range("D1:D1048576").Select() | Out-Null
